$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the organization website cell (B10): www.stat.kg -> www.stat.gov.kg,
# turn it into a real hyperlink like the neighbouring reference cells (B8, B26).
# Copy the existing hyperlink-cell formatting (from B26) first so the result
# reuses the same "Hyperlink" style used elsewhere in the sheet.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Hyperlinks.Add($ws.Range("B10"), "http://www.stat.gov.kg/", "", "", "www.stat.gov.kg")
$ws.Range("B26").Copy()
$ws.Range("B10").PasteSpecial(-4122)

# Move the current selection to B4 (matches the saved cursor position).
$ws.Range("B4").Select()
